# "Generate Report for Handback" - populate the handback columns (Latest
# Target File / Latest Handback File / Latest Handback DateTime) for the
# zh-cn and de-de localization-status sheets, flip the Status column text
# from "Ready for handoff" to "Handed back: in sync with en-US", and widen
# the columns that now hold longer file names.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Status column: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# --- Column widths: widen columns that now hold long handback file names ---
# NOTE: this engine's Columns.ColumnWidth setter stores widths in 1/6-character
# (pixel) increments (stored_width = ColumnWidth + 5/6, rounded to the nearest
# 1/6). The target stored widths below are 29.9777047293527 (rounds to 30) and
# 40 (already exact), so we feed it the ColumnWidth that lands on those grid
# points.
$wideColWidth = 30 - 0.8333333333333334       # -> stored width 30 (closest grid point to 29.9777047293527)
$maxColWidth = 40 - 0.8333333333333334        # -> stored width 40 (exact)

$overview.Columns.Item(5).ColumnWidth = $wideColWidth
$overview.Columns.Item(6).ColumnWidth = $wideColWidth

$zhcn.Columns.Item(3).ColumnWidth = $wideColWidth
$zhcn.Columns.Item(9).ColumnWidth = $maxColWidth
$zhcn.Columns.Item(10).ColumnWidth = $maxColWidth

$dede.Columns.Item(3).ColumnWidth = $wideColWidth
$dede.Columns.Item(9).ColumnWidth = $maxColWidth
$dede.Columns.Item(10).ColumnWidth = $maxColWidth

# --- zh-cn report: rows 2 & 3 (Latest Target File / Latest Handback File / Latest Handback DateTime) ---
$zhcnRows = @(
    @{ Row = 2; Name = "ba842621-2f83-45f8-a1f8-65c4195fc944.md"; TargetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f526013f552ef2bbb2b1772eaa8732b805e6c377/e2e/ba842621-2f83-45f8-a1f8-65c4195fc944.md"; HandbackFile = "ba842621-2f83-45f8-a1f8-65c4195fc944.bc428573dbb426425841c7bb35409389b1670123.zh-cn.xlf" },
    @{ Row = 3; Name = "e9195a20-65c0-4e85-951e-e7e97558922a.md"; TargetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f526013f552ef2bbb2b1772eaa8732b805e6c377/e2e/e9195a20-65c0-4e85-951e-e7e97558922a.md"; HandbackFile = "e9195a20-65c0-4e85-951e-e7e97558922a.f312c8649d8eece577746563370ac64149feb282.zh-cn.xlf" }
)

foreach ($r in $zhcnRows) {
    $iCell = $zhcn.Cells.Item($r.Row, 9)
    $iCell.Value = $r.Name
    $zhcn.Hyperlinks.Add($iCell, $r.TargetUrl, "", "", $r.Name)
    $zhcn.Cells.Item($r.Row, 10).Value = $r.HandbackFile
}
$zhcn.Range("K2").Value = "2016-08-21 23:08:43"
$zhcn.Range("K3").Value = "2016-08-21 23:08:43"

# --- de-de report: rows 2 & 3 ---
$dedeRows = @(
    @{ Row = 2; Name = "ba842621-2f83-45f8-a1f8-65c4195fc944.md"; TargetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f526013f552ef2bbb2b1772eaa8732b805e6c377/e2e/ba842621-2f83-45f8-a1f8-65c4195fc944.md"; HandbackFile = "ba842621-2f83-45f8-a1f8-65c4195fc944.bc428573dbb426425841c7bb35409389b1670123.de-de.xlf" },
    @{ Row = 3; Name = "e9195a20-65c0-4e85-951e-e7e97558922a.md"; TargetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f526013f552ef2bbb2b1772eaa8732b805e6c377/e2e/e9195a20-65c0-4e85-951e-e7e97558922a.md"; HandbackFile = "e9195a20-65c0-4e85-951e-e7e97558922a.f312c8649d8eece577746563370ac64149feb282.de-de.xlf" }
)

foreach ($r in $dedeRows) {
    $iCell = $dede.Cells.Item($r.Row, 9)
    $iCell.Value = $r.Name
    $dede.Hyperlinks.Add($iCell, $r.TargetUrl, "", "", $r.Name)
    $dede.Cells.Item($r.Row, 10).Value = $r.HandbackFile
}
$dede.Range("K2").Value = "2016-08-21 23:08:50"
$dede.Range("K3").Value = "2016-08-21 23:08:50"

Write-Output "edits applied"
